$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (shared by Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: record the generated handback target/handback file + datetime
# ---------------------------------------------------------------------------
$mdName = "4cc044c5-93a8-4f3c-9be7-ec7391096bff.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/4cc044c5-93a8-4f3c-9be7-ec7391096bff.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdName)
$wsZhCn.Range("J2").Value = "4cc044c5-93a8-4f3c-9be7-ec7391096bff.8b3f76dc52855dfb4c037a99b9d7d787ea63238a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-26 15:12:43"

# ---------------------------------------------------------------------------
# de-de sheet: record the generated handback target/handback file + datetime
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdName)
$wsDeDe.Range("J2").Value = "4cc044c5-93a8-4f3c-9be7-ec7391096bff.8b3f76dc52855dfb4c037a99b9d7d787ea63238a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-26 15:12:50"

# ---------------------------------------------------------------------------
# Column widths widened to fit the longer status/report text & file names
# ---------------------------------------------------------------------------
$wideStatusWidth = 29.166666666666668   # renders as ~29.98 character units
$wideFileWidth   = 39.166666666666664   # renders as 40 character units

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideFileWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideFileWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth
